$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in A6 (extends used range / dimension to A1:A6)
$ws.Range("A6").Value = 22222222

# Move the active selection to H7, matching the target selection in the diff
$ws.Range("H7").Select()
